$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3938.6155
$ws.Range("I74").Value = 3900.25
$ws.Range("K74").Value = 3900.25
$ws.Range("M74").Value = -2964.25
$ws.Range("H77").Value = 3938.6155
$ws.Range("I77").Value = 3900.25
$ws.Range("K77").Value = 19501.25
$ws.Range("M77").Value = -14821.25
$ws.Range("H92").Value = 2180.9
$ws.Range("I92").Value = 950.75
$ws.Range("J92").Value = 3001
$ws.Range("K92").Value = 950.75
$ws.Range("L92").Value = 3001
$ws.Range("M92").Value = 297.25
$ws.Range("N92").Value = -5497
$ws.Range("H107").Value = 816.4
$ws.Range("I107").Value = 489.75
$ws.Range("J107").Value = 1034.1666
$ws.Range("K107").Value = 489.75
$ws.Range("L107").Value = 1034.1666
$ws.Range("M107").Value = 1430.25
$ws.Range("N107").Value = -4874.1666
$ws.Range("H137").Value = 2502160.2
$ws.Range("I137").Value = 3126781.2
$ws.Range("J137").Value = 3676.125
$ws.Range("K137").Value = 9380343.600000001
$ws.Range("L137").Value = 11028.375
$ws.Range("M137").Value = -9377793.600000001
$ws.Range("N137").Value = -16128.375
$ws.Range("H138").Value = 2647789.5
$ws.Range("J138").Value = 2875771.2
$ws.Range("L138").Value = 8627313.600000001
$ws.Range("N138").Value = -8637593.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 143143740
$ws.Range("I61").Value = 250250500
$ws.Range("J61").Value = 334740.66
$ws.Range("K61").Value = 250250500
$ws.Range("L61").Value = 334740.66
$ws.Range("M61").Value = -250250288
$ws.Range("N61").Value = -335164.66
$ws.Range("H122").Value = 1558.1818
$ws.Range("I122").Value = 1656.8889
$ws.Range("J122").Value = 1114
$ws.Range("K122").Value = 4970.6667
$ws.Range("L122").Value = 3342
$ws.Range("M122").Value = -2520.6667
$ws.Range("N122").Value = -8242
$ws.Range("H132").Value = 56272.758
$ws.Range("I132").Value = 40295.383
$ws.Range("K132").Value = 120886.149
$ws.Range("M132").Value = -118356.149
$ws.Range("H136").Value = 143143740
$ws.Range("I136").Value = 250250500
$ws.Range("J136").Value = 334740.66
$ws.Range("K136").Value = 750751500
$ws.Range("L136").Value = 1004221.98
$ws.Range("M136").Value = -750748950
$ws.Range("N136").Value = -1009321.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 567.63635
$ws.Range("I80").Value = 119.5
$ws.Range("J80").Value = 667.2222
$ws.Range("K80").Value = 119.5
$ws.Range("L80").Value = 667.2222
$ws.Range("M80").Value = 878.5
$ws.Range("N80").Value = -2663.2222
$ws.Range("H83").Value = 567.63635
$ws.Range("I83").Value = 119.5
$ws.Range("J83").Value = 667.2222
$ws.Range("K83").Value = 597.5
$ws.Range("L83").Value = 3336.111
$ws.Range("M83").Value = 4394.5
$ws.Range("N83").Value = -13320.111
$ws.Range("H86").Value = 12904.667
$ws.Range("I86").Value = 19068.77
$ws.Range("J86").Value = 2888
$ws.Range("K86").Value = 19068.77
$ws.Range("L86").Value = 2888
$ws.Range("M86").Value = -17945.77
$ws.Range("N86").Value = -5134
$ws.Range("H89").Value = 12904.667
$ws.Range("I89").Value = 19068.77
$ws.Range("J89").Value = 2888
$ws.Range("K89").Value = 95343.85000000001
$ws.Range("L89").Value = 14440
$ws.Range("M89").Value = -89727.85000000001
$ws.Range("N89").Value = -25672
$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 2000
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2902
$ws.Range("H107").Value = 2375.125
$ws.Range("I107").Value = 1747.5
$ws.Range("J107").Value = 2584.3333
$ws.Range("K107").Value = 1747.5
$ws.Range("L107").Value = 2584.3333
$ws.Range("M107").Value = 172.5
$ws.Range("N107").Value = -6424.3333
$ws.Range("H134").Value = 2005.5
$ws.Range("I134").Value = 1928.1818
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 5784.5454
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -3249.5454
$ws.Range("N134").Value = -11370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2216.415
$ws.Range("I31").Value = 1106.7179
$ws.Range("K31").Value = 1106.7179
$ws.Range("M31").Value = -811.7179000000001
$ws.Range("H34").Value = 2216.415
$ws.Range("I34").Value = 1106.7179
$ws.Range("K34").Value = 1106.7179
$ws.Range("M34").Value = -904.7179000000001
$ws.Range("H99").Value = 4762.524
$ws.Range("I99").Value = 5001.091
$ws.Range("J99").Value = 4500.1
$ws.Range("K99").Value = 5001.091
$ws.Range("L99").Value = 4500.1
$ws.Range("M99").Value = -3503.091
$ws.Range("N99").Value = -7496.1
$ws.Range("H107").Value = 1416.6666
$ws.Range("I107").Value = 625
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 625
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 1295
$ws.Range("N107").Value = -6840
$ws.Range("H126").Value = 4762.524
$ws.Range("I126").Value = 5001.091
$ws.Range("J126").Value = 4500.1
$ws.Range("K126").Value = 15003.273
$ws.Range("L126").Value = 13500.3
$ws.Range("M126").Value = -12533.273
$ws.Range("N126").Value = -18440.3
$ws.Range("H127").Value = 31842.105
$ws.Range("J127").Value = 31842.105
$ws.Range("L127").Value = 31842.105
$ws.Range("N127").Value = -41762.105
$ws.Range("H132").Value = 46022.957
$ws.Range("I132").Value = 2963
$ws.Range("J132").Value = 79146
$ws.Range("K132").Value = 8889
$ws.Range("L132").Value = 237438
$ws.Range("M132").Value = -6359
$ws.Range("N132").Value = -242498
$ws.Range("H134").Value = 53605.43
$ws.Range("I134").Value = 2362.5
$ws.Range("J134").Value = 217582.8
$ws.Range("K134").Value = 7087.5
$ws.Range("L134").Value = 652748.3999999999
$ws.Range("M134").Value = -4552.5
$ws.Range("N134").Value = -657818.3999999999
$ws.Range("H137").Value = 33450
$ws.Range("J137").Value = 33450
$ws.Range("L137").Value = 33450
$ws.Range("N137").Value = -43650

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1400
$ws.Range("I68").Value = 769.2308
$ws.Range("J68").Value = 5500
$ws.Range("K68").Value = 2307.6924
$ws.Range("L68").Value = 16500
$ws.Range("M68").Value = -1496.6924
$ws.Range("N68").Value = -18122
$ws.Range("H71").Value = 1400
$ws.Range("I71").Value = 769.2308
$ws.Range("J71").Value = 5500
$ws.Range("K71").Value = 6923.077200000001
$ws.Range("L71").Value = 49500
$ws.Range("M71").Value = -2867.077200000001
$ws.Range("N71").Value = -57612
$ws.Range("H80").Value = 4216.6665
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4216.6665
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 12649.9995
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -14521.9995
$ws.Range("H83").Value = 4216.6665
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4216.6665
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 37949.9985
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -47309.9985
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("H92").Value = 1077.2222
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 1465.8334
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 4397.5002
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -6893.5002
$ws.Range("H131").Value = 705
$ws.Range("I131").Value = 590
$ws.Range("J131").Value = 992.5
$ws.Range("K131").Value = 1770
$ws.Range("L131").Value = 2977.5
$ws.Range("M131").Value = 3270
$ws.Range("N131").Value = -13057.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4826.7075
$ws.Range("I70").Value = 4803.9585
$ws.Range("J70").Value = 4858.8237
$ws.Range("K70").Value = 4803.9585
$ws.Range("L70").Value = 4858.8237
$ws.Range("M70").Value = -4533.9585
$ws.Range("N70").Value = -5398.8237
$ws.Range("H73").Value = 4826.7075
$ws.Range("I73").Value = 4803.9585
$ws.Range("J73").Value = 4858.8237
$ws.Range("K73").Value = 4803.9585
$ws.Range("L73").Value = 4858.8237
$ws.Range("M73").Value = -3867.9585
$ws.Range("N73").Value = -6730.8237

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 252.28572
$ws.Range("I55").Value = 230.76471
$ws.Range("J55").Value = 301.06668
$ws.Range("K55").Value = 230.76471
$ws.Range("L55").Value = 301.06668
$ws.Range("M55").Value = -57.76471000000001
$ws.Range("N55").Value = -647.06668
$ws.Range("H93").Value = 1338.3077
$ws.Range("I93").Value = 965.55554
$ws.Range("J93").Value = 2177
$ws.Range("K93").Value = 965.55554
$ws.Range("L93").Value = 2177
$ws.Range("M93").Value = 282.44446
$ws.Range("N93").Value = -4673

